$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "47.326.29"
$ws.Range("E2").Value = "  +2.47%  "

# Row 3
$ws.Range("D3").Value = "2.502.42"
$ws.Range("E3").Value = "  +2.04%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").Value = "323.41"
$ws.Range("E5").Value = "  +0.39%  "

# Row 6
$ws.Range("D6").Value = "108.78"
$ws.Range("E6").Value = "  +3.52%  "

# Row 8
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  -0.16%  "

# Row 10
$ws.Range("D10").Value = "39.11"
$ws.Range("E10").Value = "  +8.64%  "

# Row 11
$ws.Range("E11").Value = "  +0.76%  "

# Row 12
$ws.Range("E12").Value = "  +0.67%  "

# Row 13
$ws.Range("D13").Value = "18.36"
$ws.Range("E13").Value = "  +0.39%  "

# Row 14
$ws.Range("D14").Value = "7.19"
$ws.Range("E14").Value = "  +1.57%  "

# Row 15
$ws.Range("D15").Value = "2.893.36"
$ws.Range("E15").Value = "  +2.12%  "

# Row 16
$ws.Range("D16").Value = "2.502.44"
$ws.Range("E16").Value = "  +2.34%  "

# Row 18
$ws.Range("D18").Value = "47.253.79"
$ws.Range("E18").Value = "  +2.56%  "

# Row 19
$ws.Range("E19").Value = "  +2.00%  "

# Row 20
$ws.Range("D20").Value = "6.65"
$ws.Range("E20").Value = "  +3.51%  "

# Row 21
$ws.Range("E21").Value = "  +0.65%  "

# Row 22
$ws.Range("D22").Value = "2.71"
$ws.Range("E22").Value = "  +12.95%  "

# Row 23
$ws.Range("D23").Value = "70.52"
$ws.Range("E23").Value = "  -0.62%  "

# Row 24
$ws.Range("D24").Value = "247.85"
$ws.Range("E24").Value = "  +0.12%  "

# Row 25
$ws.Range("E25").Value = "  +3.30%  "

# Row 26
$ws.Range("D26").Value = "26.07"
$ws.Range("E26").Value = "  +0.56%  "

# Row 27
$ws.Range("E27").Value = "  -0.01%  "

# Row 28
$ws.Range("E28").Value = "  +0.34%  "

# Row 29
$ws.Range("E29").Value = "  +3.85%  "

# Row 30
$ws.Range("D30").Value = "35.37"
$ws.Range("E30").Value = "  +2.62%  "

# Row 31
$ws.Range("E31").Value = "  +7.45%  "

# Row 32
$ws.Range("D32").Value = "49.85"
$ws.Range("E32").Value = "  +1.06%  "

# Row 33
$ws.Range("D33").Value = "20.01"
$ws.Range("E33").Value = "  +1.16%  "

# Row 34
$ws.Range("E34").Value = "  +1.80%  "

# Row 35
$ws.Range("E35").Value = "  +3.29%  "

# Row 36
$ws.Range("E36").Value = "  +0.22%  "

# Row 37
$ws.Range("E37").Value = "  +4.47%  "

# Row 38
$ws.Range("D38").Value = "'4.70"
$ws.Range("E38").Value = "  +3.49%  "

# Row 39
$ws.Range("E39").Value = "  +1.15%  "

# Row 40
$ws.Range("E40").Value = "  +1.11%  "

# Row 41
$ws.Range("E41").Value = "  +0.51%  "

# Row 42
$ws.Range("D42").Value = "120.96"
$ws.Range("E42").Value = "  -5.33%  "

# Row 43
$ws.Range("D43").Value = "21.32"
$ws.Range("E43").Value = "  +1.93%  "

# Row 44
$ws.Range("E44").Value = "  +2.03%  "

# Row 45
$ws.Range("D45").Value = "1.991.08"
$ws.Range("E45").Value = "  +1.03%  "

# Row 46
$ws.Range("D46").Value = "3.06"
$ws.Range("E46").Value = "  +3.14%  "

# Row 47
$ws.Range("D47").Value = "2.06"
$ws.Range("E47").Value = "  -1.49%  "

# Row 48
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "1.78"
$ws.Range("E48").Value = "  -3.90%  "

# Row 49
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "9.09"
$ws.Range("E49").Value = "  +0.19%  "

# Row 50
$ws.Range("D50").Value = "5.22"
$ws.Range("E50").Value = "  +3.22%  "

# Row 51
$ws.Range("D51").Value = "56.46"
$ws.Range("E51").Value = "  +3.69%  "
